# canasta_diaria.xlsx — "renta imputada para casa cedida o prestada"
#
# Upstream methodology change: households living in housing that is ceded
# or lent (casa cedida o prestada) now get an imputed rent, which shifts
# the survey population used for the daily food-basket (canasta diaria)
# calculation from 3005 to 3009 observations. That ripples through every
# row: population (D), cantidad_h (B), cantidad_ajustada (C) and cal_intake
# (F) are all recomputed; cal (E) is untouched except where the underlying
# item's sort position moved (rows 7/8, 15/16, 18/19), in which case the
# item label (column A) and its "cal" figure travel together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: re-label the rows whose item swapped position -------------
$ws.Range("A7").Value  = "Pastas alimenticias"
$ws.Range("A8").Value  = "Yuca"
$ws.Range("A15").Value = "Frijoles"
$ws.Range("A16").Value = "Platanos"
$ws.Range("A18").Value = "Caraotas"
$ws.Range("A19").Value = "Leche en polvo, completa o descremada"

# --- Columns B-F: recomputed cantidad_h / cantidad_ajustada / population /
#     cal / cal_intake values for every data row (rows 2-29) --------------

$ws.Range("B2").Value = 94.976976093559728
$ws.Range("C2").Value = 117.96357444086199
$ws.Range("D2").Value = 3009
$ws.Range("E2").Value = 383
$ws.Range("F2").Value = 451.80047607421875

$ws.Range("B3").Value = 88.071501905872964
$ws.Range("C3").Value = 109.386817814227
$ws.Range("D3").Value = 3009
$ws.Range("E3").Value = 345
$ws.Range("F3").Value = 377.384521484375

$ws.Range("B4").Value = 19.665394595351557
$ws.Range("C4").Value = 24.424869165146379
$ws.Range("D4").Value = 3009
$ws.Range("E4").Value = 900
$ws.Range("F4").Value = 219.82382202148438

$ws.Range("B5").Value = 29.394199233785663
$ws.Range("C5").Value = 36.50826673041788
$ws.Range("D5").Value = 3009
$ws.Range("E5").Value = 393.5
$ws.Range("F5").Value = 143.6600341796875

$ws.Range("B6").Value = 20.209372464348288
$ws.Range("C6").Value = 25.100502154881791
$ws.Range("D6").Value = 3009
$ws.Range("E6").Value = 368.5
$ws.Range("F6").Value = 92.495353698730469

$ws.Range("B7").Value = 45.010208802747741
$ws.Range("C7").Value = 55.903707537868165
$ws.Range("D7").Value = 3009
$ws.Range("E7").Value = 137.5
$ws.Range("F7").Value = 76.867599487304688

$ws.Range("B8").Value = 33.815221819460845
$ws.Range("C8").Value = 41.999278094214674
$ws.Range("D8").Value = 3009
$ws.Range("E8").Value = 182.33332824707031
$ws.Range("F8").Value = 76.578681945800781

$ws.Range("B9").Value = 28.417130347936801
$ws.Range("C9").Value = 35.294724975007838
$ws.Range("D9").Value = 3009
$ws.Range("E9").Value = 196.5
$ws.Range("F9").Value = 69.354133605957031

$ws.Range("B10").Value = 14.018136161076503
$ws.Range("C10").Value = 17.410845444581959
$ws.Range("D10").Value = 3009
$ws.Range("E10").Value = 355
$ws.Range("F10").Value = 61.808502197265625

$ws.Range("B11").Value = 17.566349087224207
$ws.Range("C11").Value = 21.817806737056085
$ws.Range("D11").Value = 3009
$ws.Range("E11").Value = 254.5
$ws.Range("F11").Value = 55.526317596435547

$ws.Range("B12").Value = 5.620282138927152
$ws.Range("C12").Value = 6.9805187861283144
$ws.Range("D12").Value = 3009
$ws.Range("E12").Value = 584
$ws.Range("F12").Value = 40.766231536865234

$ws.Range("B13").Value = 28.7704038673875
$ws.Range("C13").Value = 35.733498839216402
$ws.Range("D13").Value = 3009
$ws.Range("E13").Value = 113.375
$ws.Range("F13").Value = 40.512855529785156

$ws.Range("B14").Value = 18.268528258614303
$ws.Range("C14").Value = 22.689929412583172
$ws.Range("D14").Value = 3009
$ws.Range("E14").Value = 174
$ws.Range("F14").Value = 39.480476379394531

$ws.Range("B15").Value = 7.4372123992098409
$ws.Range("C15").Value = 9.2371876337235843
$ws.Range("D15").Value = 3009
$ws.Range("E15").Value = 405.84616088867188
$ws.Range("F15").Value = 37.48876953125

$ws.Range("B16").Value = 16.781085621918425
$ws.Range("C16").Value = 20.842491949860459
$ws.Range("D16").Value = 3009
$ws.Range("E16").Value = 164.85714721679688
$ws.Range("F16").Value = 34.360336303710938

$ws.Range("B17").Value = 16.742273180458625
$ws.Range("C17").Value = 20.794285805291139
$ws.Range("D17").Value = 3009
$ws.Range("E17").Value = 145
$ws.Range("F17").Value = 30.151714324951172

$ws.Range("B18").Value = 16.975265168240394
$ws.Range("C18").Value = 21.083666997004997
$ws.Range("D18").Value = 3009
$ws.Range("E18").Value = 135.11111450195313
$ws.Range("F18").Value = 28.486377716064453

$ws.Range("B19").Value = 5.3506149702113346
$ws.Range("C19").Value = 6.6455859827337651
$ws.Range("D19").Value = 3009
$ws.Range("E19").Value = 428.5
$ws.Range("F19").Value = 28.476335525512695

$ws.Range("B20").Value = 26.817168198696756
$ws.Range("C20").Value = 33.307535595575338
$ws.Range("D20").Value = 3009
$ws.Range("E20").Value = 85
$ws.Range("F20").Value = 28.311405181884766

$ws.Range("B21").Value = 9.8233872607285964
$ws.Range("C21").Value = 12.200871343666551
$ws.Range("D21").Value = 3009
$ws.Range("E21").Value = 183.25
$ws.Range("F21").Value = 22.358097076416016

$ws.Range("B22").Value = 13.281584191060851
$ws.Range("C22").Value = 16.496030918424331
$ws.Range("D22").Value = 3009
$ws.Range("E22").Value = 122.46154022216797
$ws.Range("F22").Value = 20.2012939453125

$ws.Range("B23").Value = 13.179509357821312
$ws.Range("C23").Value = 16.369251651199761
$ws.Range("D23").Value = 3009
$ws.Range("E23").Value = 40
$ws.Range("F23").Value = 6.5477008819580078

$ws.Range("B24").Value = 1.5112757105968371
$ws.Range("C24").Value = 1.8770389899180633
$ws.Range("D24").Value = 3009
$ws.Range("E24").Value = 284.66665649414063
$ws.Range("F24").Value = 5.3433041572570801

$ws.Range("B25").Value = 6.6915920387672481
$ws.Range("C25").Value = 8.3111101554928926
$ws.Range("D25").Value = 3009
$ws.Range("E25").Value = 56.5
$ws.Range("F25").Value = 4.695777416229248

$ws.Range("B26").Value = 9.9464465566313756
$ws.Range("C26").Value = 12.353713879834064
$ws.Range("D26").Value = 3009
$ws.Range("E26").Value = 35
$ws.Range("F26").Value = 4.3238000869750977

$ws.Range("B27").Value = 8.7791864503535297
$ws.Range("C27").Value = 10.903950118828391
$ws.Range("D27").Value = 3009
$ws.Range("E27").Value = 23
$ws.Range("F27").Value = 2.5079085826873779

$ws.Range("B28").Value = 11.080093315171419
$ws.Range("C28").Value = 13.761728929815677
$ws.Range("D28").Value = 3009
$ws.Range("E28").Value = 5
$ws.Range("F28").Value = 0.68808645009994507

$ws.Range("B29").Value = 17.59972508986078
$ws.Range("C29").Value = 21.859260514234933
$ws.Range("D29").Value = 3009
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
